$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 10
$ws.Range("I31").Value = 10
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 30
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 200
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H138").Value = 2975.4783
$ws.Range("I138").Value = 1484.2222
$ws.Range("J138").Value = 3934.1428
$ws.Range("K138").Value = 4452.6666
$ws.Range("L138").Value = 11802.4284
$ws.Range("M138").Value = 687.3334000000004

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2013.4419
$ws.Range("I32").Value = 1855.561
$ws.Range("J32").Value = 5250
$ws.Range("K32").Value = 1855.561
$ws.Range("L32").Value = 5250
$ws.Range("M32").Value = -1568.561
$ws.Range("H45").Value = 4100
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4100
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4100
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4854
$ws.Range("H61").Value = 2619.9333
$ws.Range("I61").Value = 1412.5
$ws.Range("J61").Value = 3999.8572
$ws.Range("K61").Value = 1412.5
$ws.Range("L61").Value = 3999.8572
$ws.Range("M61").Value = -1200.5
$ws.Range("N61").Value = -4423.8572
$ws.Range("H74").Value = 1713.5
$ws.Range("I74").Value = 1368
$ws.Range("J74").Value = 2750
$ws.Range("K74").Value = 1368
$ws.Range("L74").Value = 2750
$ws.Range("M74").Value = -494
$ws.Range("H77").Value = 1713.5
$ws.Range("I77").Value = 1368
$ws.Range("J77").Value = 2750
$ws.Range("K77").Value = 6840
$ws.Range("L77").Value = 13750
$ws.Range("M77").Value = -2472
$ws.Range("H97").Value = 2545.6667
$ws.Range("I97").Value = 2854.8
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2854.8
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -2358.8
$ws.Range("H109").Value = 75000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 75000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77774
$ws.Range("H136").Value = 2619.9333
$ws.Range("I136").Value = 1412.5
$ws.Range("J136").Value = 3999.8572
$ws.Range("K136").Value = 4237.5
$ws.Range("L136").Value = 11999.5716
$ws.Range("M136").Value = -1687.5
$ws.Range("N136").Value = -17099.5716

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 15333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 15333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 15333
$ws.Range("N92").Value = -20325
$ws.Range("H94").Value = 1955.2
$ws.Range("I94").Value = 1972
$ws.Range("J94").Value = 1888
$ws.Range("K94").Value = 1972
$ws.Range("L94").Value = 1888
$ws.Range("M94").Value = -1521
$ws.Range("N94").Value = -2790
$ws.Range("H134").Value = 1982.7142
$ws.Range("I134").Value = 1194.5714
$ws.Range("J134").Value = 3559
$ws.Range("K134").Value = 3583.7142
$ws.Range("L134").Value = 10677
$ws.Range("M134").Value = -1048.7142
$ws.Range("N134").Value = -15747

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H58").Value = 1496.4
$ws.Range("I58").Value = 1498.2222
$ws.Range("J58").Value = 1480
$ws.Range("K58").Value = 1498.2222
$ws.Range("L58").Value = 1480
$ws.Range("M58").Value = -1295.2222
$ws.Range("H105").Value = 3802
$ws.Range("I105").Value = 3502.5
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3502.5
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1755.5
$ws.Range("N105").Value = -8494
$ws.Range("H122").Value = 2325.2
$ws.Range("I122").Value = 1904
$ws.Range("J122").Value = 2957
$ws.Range("K122").Value = 5712
$ws.Range("L122").Value = 8871
$ws.Range("M122").Value = -3262
$ws.Range("N122").Value = -13771
$ws.Range("H134").Value = 1401.5714
$ws.Range("I134").Value = 1485.1666
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 4455.4998
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -1920.4998
$ws.Range("N134").Value = -7770
$ws.Range("H136").Value = 1496.4
$ws.Range("I136").Value = 1498.2222
$ws.Range("J136").Value = 1480
$ws.Range("K136").Value = 4494.6666
$ws.Range("L136").Value = 4440
$ws.Range("M136").Value = -1944.6666

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 182.33333
$ws.Range("I107").Value = 223.5
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 670.5
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1249.5
$ws.Range("N107").Value = -4140

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4248.75
$ws.Range("I113").Value = 4141.4287
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4141.4287
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1971.4287
$ws.Range("H132").Value = 2401.5715
$ws.Range("I132").Value = 1919.4
$ws.Range("J132").Value = 3607
$ws.Range("K132").Value = 5758.200000000001
$ws.Range("L132").Value = 10821
$ws.Range("M132").Value = -3228.200000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1995.8
$ws.Range("I22").Value = 1989.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1989.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1694.5
$ws.Range("H27").Value = 1995.8
$ws.Range("I27").Value = 1989.5
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1989.5
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1882.5
$ws.Range("H46").Value = 2552.5557
$ws.Range("I46").Value = 1829.0834
$ws.Range("J46").Value = 3999.5
$ws.Range("K46").Value = 1829.0834
$ws.Range("L46").Value = 3999.5
$ws.Range("M46").Value = -1641.0834
$ws.Range("N46").Value = -4375.5
$ws.Range("H55").Value = 2556.1428
$ws.Range("I55").Value = 2565.5
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 2565.5
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = -2392.5
$ws.Range("H61").Value = 2003.5714
$ws.Range("I61").Value = 2003.5714
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2003.5714
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1801.5714
$ws.Range("H113").Value = 2003.5714
$ws.Range("I113").Value = 2003.5714
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2003.5714
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 166.4286
$ws.Range("H136").Value = 809.7143
$ws.Range("I136").Value = 809.7143
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2429.1429
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 120.8571000000002

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 15019
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 15019
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 15019
$ws.Range("N28").Value = -15715
$ws.Range("H81").Value = 11479.167
$ws.Range("I81").Value = 13575
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 27150
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -26089
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 11479.167
$ws.Range("I84").Value = 13575
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 135750
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -130446
$ws.Range("N84").Value = -20608
$ws.Range("H132").Value = 2543.484
$ws.Range("I132").Value = 1342.85
$ws.Range("J132").Value = 4726.4546
$ws.Range("K132").Value = 4028.55
$ws.Range("L132").Value = 14179.3638
$ws.Range("M132").Value = -1498.55
